# Master_QBR_Template.pptx - Strategic Recommendations slide
#
# Re-lays out the "Strategic Recommendations" slide (slide 7) so each
# numbered recommendation gets a bold "title" line plus a new smaller
# "rationale" line underneath it (populated from the LLM). The three
# existing {{RECOMMENDATION_n}} placeholders are replaced by a
# {{REC_n_TITLE}} / {{REC_n_RATIONALE}} pair, and everything is shifted
# up / tightened to make room.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# The PowerPoint object model works in points (1 pt = 12700 EMU) for all
# shape geometry, but the host's internal point<->EMU round trip can be
# off by a single EMU because of floating point truncation. Nudge every
# converted value by a hair so it lands back on the exact EMU we want.
function EMU([double]$emu) {
    return ($emu / 12700.0) + 0.00005
}

# Grab stable references to every shape we need to touch *before* we
# start inserting new shapes - inserting shifts collection indices, but
# object references stay bound to the same underlying shape.
$title = $s.Shapes.Item(1)

$rect1 = $s.Shapes.Item(2)
$num1 = $s.Shapes.Item(3)
$rec1 = $s.Shapes.Item(4)

$rect2 = $s.Shapes.Item(5)
$num2 = $s.Shapes.Item(6)
$rec2 = $s.Shapes.Item(7)

$rect3 = $s.Shapes.Item(8)
$num3 = $s.Shapes.Item(9)
$rec3 = $s.Shapes.Item(10)

$rects = @($rect1, $rect2, $rect3)
$nums = @($num1, $num2, $num3)
$recs = @($rec1, $rec2, $rec3)

# Row geometry, in EMU, taken straight from the target layout.
$rowTops = @(1097280, 2011680, 2926080)
$numSize = 365760
$titleLeft = 1005840
$titleWidth = 7589520
$titleHeight = 365760
$rationaleTops = @(1481328, 2395728, 3310128)
$rationaleHeight = 457200

$titleTexts = @("{{REC_1_TITLE}}", "{{REC_2_TITLE}}", "{{REC_3_TITLE}}")
$rationaleTexts = @("{{REC_1_RATIONALE}}", "{{REC_2_RATIONALE}}", "{{REC_3_RATIONALE}}")
$rationaleNames = @("TextBox 5", "TextBox 9", "TextBox 13")

# ---------------------------------------------------------------------
# Title banner ("Strategic Recommendations") - shrink + move up a touch
# ---------------------------------------------------------------------
$title.Left = EMU 457200
$title.Top = EMU 274320
$title.Width = EMU 8229600
$title.Height = EMU 640080
$title.TextFrame.TextRange.Font.Size = 36

# ---------------------------------------------------------------------
# Each of the 3 recommendation rows
# ---------------------------------------------------------------------
for ($row = 0; $row -lt 3; $row++) {
    $top = $rowTops[$row]

    # --- number background rectangle -> smaller 28.8 x 28.8 pt square ---
    $rect = $rects[$row]
    $rect.Left = EMU 457200
    $rect.Top = EMU $top
    $rect.Width = EMU $numSize
    $rect.Height = EMU $numSize

    # --- number textbox (the "1"/"2"/"3") -> same smaller square, smaller font ---
    $num = $nums[$row]
    $num.Left = EMU 457200
    $num.Top = EMU $top
    $num.Width = EMU $numSize
    $num.Height = EMU $numSize
    $num.TextFrame.TextRange.Font.Size = 16

    # --- recommendation textbox becomes the bold "title" line ---
    $rec = $recs[$row]
    $rec.Left = EMU $titleLeft
    $rec.Top = EMU $top
    $rec.Width = EMU $titleWidth
    $rec.TextFrame.TextRange.Text = $titleTexts[$row]
    $rec.TextFrame.TextRange.Font.Size = 14
    $rec.TextFrame.TextRange.Font.Bold = -1
    $rec.TextFrame.TextRange.Font.Color.RGB = 0x8A5C2E
    # Set height last - the shape has AutoSize/spAutoFit enabled, so
    # earlier font-size changes recompute it; pin it to the template value.
    $rec.Height = EMU $titleHeight

    # --- new rationale textbox, inserted right after the title line ---
    $rat = $s.Shapes.AddTextbox(1, (EMU $titleLeft), (EMU $rationaleTops[$row]), (EMU $titleWidth), (EMU $rationaleHeight))
    $rat.Name = $rationaleNames[$row]
    $rat.Fill.Visible = 0
    $rat.TextFrame.WordWrap = -1
    $rat.TextFrame.AutoSize = 1
    $rat.TextFrame.TextRange.Text = $rationaleTexts[$row]
    $rat.TextFrame.TextRange.Font.Size = 11
    $rat.TextFrame.TextRange.Font.Color.RGB = 0x68554A
    # Force the height/position back to the template's fixed values
    # (AutoSize recalculates height as text/font are applied above).
    $rat.Left = EMU $titleLeft
    $rat.Top = EMU $rationaleTops[$row]
    $rat.Width = EMU $titleWidth
    $rat.Height = EMU $rationaleHeight

    # Move the newly-appended shape so it sits immediately after the
    # recommendation/title textbox it belongs to, instead of at the
    # very end of the z-order.
    $targetPos = $rec.ZOrderPosition + 1
    while ($rat.ZOrderPosition -gt $targetPos) {
        $rat.ZOrder(3)
    }
}
